# Generate Report for Handoff
# Updates the Overview, zh-cn and de-de sheets with a fresh handoff report for
# the rows whose handoff just ran again (076bdd39, 0d1f501c, 20d780be,
# a85cd41e, d0850586, fbde40c4):
#  - Latest HO Xliff Generate Date (Overview col G) moves from 08:21:40 -> 08:21:57
#  - Latest Handoff Datetime (de-de col H) happens to carry the same original
#    timestamp as the Overview column, so it also becomes 08:21:57
#  - Latest Handoff Datetime (zh-cn col H) moves from 08:21:35 -> 08:21:51
#  - Priority (zh-cn/de-de col E) for the matching rows is now marked "ht"

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 12, 13, 14)

$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-08-27 08:21:57"
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
    $wsZhCn.Cells.Item($r, 8).Value = "2016-08-27 08:21:51"
}

$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
    $wsDeDe.Cells.Item($r, 8).Value = "2016-08-27 08:21:57"
}
